# Updated cryptos list with GitHub Actions: refresh Price (column D) and
# Volume(1h) (column E) figures for each coin row. Values that look like
# plain numbers are entered with a leading apostrophe so Excel keeps them
# as text (matching the original inlineStr cells) instead of coercing them
# into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.146.70"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "2.509.15"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'592.68"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'175.60"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.516"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "2.506.73"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "'0.151"
$ws.Range("E10").Value = "  +5.77%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "'5.00"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "'0.336"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").Value = "2.951.85"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "'25.76"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "68.959.23"
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("D17").Value = "'0.0000173"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "2.505.99"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "'361.72"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").Value = "'7.52"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "'10.93"
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'70.14"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "'4.17"
$ws.Range("E25").Value = "  -3.11%  "
$ws.Range("D26").Value = "'8.93"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").Value = "'1.66"
$ws.Range("E27").Value = "  -7.19%  "
$ws.Range("D28").Value = "2.626.52"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "'507.22"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").Value = "0.0₃0879"
$ws.Range("E31").Value = "  -3.99%  "
$ws.Range("D32").Value = "'7.70"
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").Value = "'1.77"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").Value = "'1.21"
$ws.Range("E34").Value = "  -4.91%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'162.30"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "'0.119"
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("D38").Value = "'18.67"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "'18.67"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D41").Value = "'1.31"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("D42").Value = "'1.70"
$ws.Range("E42").Value = "  -2.46%  "
$ws.Range("D43").Value = "'4.74"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D44").Value = "'0.318"
$ws.Range("E44").Value = "  -3.76%  "
$ws.Range("D45").Value = "'2.31"
$ws.Range("E45").Value = "  -4.76%  "
$ws.Range("D46").Value = "'149.47"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").Value = "'3.54"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "'0.511"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "'0.0735"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "'1.56"
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("D51").Value = "'0.577"
$ws.Range("E51").Value = "  -1.82%  "
